$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 187 ("「間違いで賢くなる、痛みで強くなる」" entry),
# which shifts all subsequent rows up by one.
$ws.Rows.Item(187).Delete()
